$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = "Maha Al Ameri"
$wsSummary.Range("B4").Value = 1869.23
$wsSummary.Range("B6").Value = 388526
$wsSummary.Range("B7").Value = 272339
$wsSummary.Range("B8").Value = 116187
$wsSummary.Range("B9").Value = 1.43

# --- Sheet 2: Assets ---
$wsAssets = $wb.Worksheets.Item("Assets")
$wsAssets.Range("B2").Value = "Luxury Car"
$wsAssets.Range("C2").Value = 385481
$wsAssets.Range("C3").Value = 3045
$wsAssets.Range("C4").Value = 388526

# --- Sheet 3: Liabilities ---
$wsLiabilities = $wb.Worksheets.Item("Liabilities")
$wsLiabilities.Range("C2").Value = 231289
$wsLiabilities.Range("D2").Value = 6425
$wsLiabilities.Range("E2").Value = 3
$wsLiabilities.Range("C3").Value = 41050
$wsLiabilities.Range("D3").Value = 2052
$wsLiabilities.Range("C4").Value = 272339
